$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 172.0013572692871
$ws.Range("C2").Value = 60.86019637923512
$ws.Range("D2").Value = 1.982184219360352
$ws.Range("E2").Value = 3.05625336877818

$ws.Range("B3").Value = 458.8654263019562
$ws.Range("C3").Value = 142.3697241792931
$ws.Range("D3").Value = 4.44204683303833
$ws.Range("E3").Value = 1.941639971233113

$ws.Range("B4").Value = 986.6406661510467
$ws.Range("C4").Value = 152.3429324197695
$ws.Range("D4").Value = 7.753339338302612
$ws.Range("E4").Value = 0.5416233107815548

$ws.Range("B5").Value = 428.5647017478943
$ws.Range("C5").Value = 2.663015539081247
$ws.Range("D5").Value = 6.445478677749634
$ws.Range("E5").Value = 0.5027006889482671

$ws.Range("B6").Value = 848.4674672603608
$ws.Range("C6").Value = 1.775739370772414
$ws.Range("D6").Value = 6.564799785614014
$ws.Range("E6").Value = 0.5584758047761286

$ws.Range("B7").Value = 1682.078295707703
$ws.Range("C7").Value = 3.637935349281217
$ws.Range("D7").Value = 7.910986375808716
$ws.Range("E7").Value = 0.3462165518427878

$ws.Range("B8").Value = 569.8244668006897
$ws.Range("C8").Value = 5.615792983125055
$ws.Range("D8").Value = 6.289184713363648
$ws.Range("E8").Value = 0.5953469587908193

$ws.Range("B9").Value = 1108.669444561005
$ws.Range("C9").Value = 2.805161702153816
$ws.Range("D9").Value = 7.41806640625
$ws.Range("E9").Value = 0.5370743424808685

$ws.Range("B10").Value = 2203.644894075393
$ws.Range("C10").Value = 9.515119678493532
$ws.Range("D10").Value = 8.372460746765137
$ws.Range("E10").Value = 0.445031818255665

$ws.Range("B11").Value = 410.3941568851471
$ws.Range("C11").Value = 6.297791436250154
$ws.Range("D11").Value = 7.739963722229004
$ws.Range("E11").Value = 0.5563292532764453

$ws.Range("B12").Value = 799.442493391037
$ws.Range("C12").Value = 14.77971849925392
$ws.Range("D12").Value = 8.079037237167359
$ws.Range("E12").Value = 0.4333297754536677

$ws.Range("B13").Value = 1587.551083898544
$ws.Range("C13").Value = 21.70231092178417
$ws.Range("D13").Value = 8.223393535614013
$ws.Range("E13").Value = 0.7165630616114224

$ws.Range("B14").Value = 591.9313493728638
$ws.Range("C14").Value = 11.18070723316746
$ws.Range("D14").Value = 8.023996877670289
$ws.Range("E14").Value = 0.4121866517894509

$ws.Range("B15").Value = 1167.719575881958
$ws.Range("C15").Value = 16.98554940408184
$ws.Range("D15").Value = 7.799665212631226
$ws.Range("E15").Value = 0.5862128792063116

$ws.Range("B16").Value = 2280.842783212662
$ws.Range("C16").Value = 31.65679145237837
$ws.Range("D16").Value = 8.221553230285645
$ws.Range("E16").Value = 0.5925313107772046

$ws.Range("B17").Value = 753.746323299408
$ws.Range("C17").Value = 19.70052129207456
$ws.Range("D17").Value = 8.324148750305175
$ws.Range("E17").Value = 0.7320789850106061

$ws.Range("B18").Value = 1529.178658294678
$ws.Range("C18").Value = 11.68163558137355
$ws.Range("D18").Value = 7.911117362976074
$ws.Range("E18").Value = 0.6596235920012946

$ws.Range("B19").Value = 3053.358260679245
$ws.Range("C19").Value = 42.70827660953847
$ws.Range("D19").Value = 8.505347776412965
$ws.Range("E19").Value = 0.5510282715553173

$ws.Range("B20").Value = 430.0886552333832
$ws.Range("C20").Value = 8.157954851886419
$ws.Range("D20").Value = 7.850925874710083
$ws.Range("E20").Value = 0.4662844979405372

$ws.Range("B21").Value = 850.224909734726
$ws.Range("C21").Value = 15.77346689464976
$ws.Range("D21").Value = 7.675096750259399
$ws.Range("E21").Value = 0.5928059403103265

$ws.Range("B22").Value = 1655.931116008759
$ws.Range("C22").Value = 34.63639089138093
$ws.Range("D22").Value = 8.676648759841919
$ws.Range("E22").Value = 0.1655224430693714

$ws.Range("B23").Value = 616.8991370677948
$ws.Range("C23").Value = 13.83232493307794
$ws.Range("D23").Value = 7.785571622848511
$ws.Range("E23").Value = 0.3129624034840556

$ws.Range("B24").Value = 1207.195433235168
$ws.Range("C24").Value = 27.17502444467787
$ws.Range("D24").Value = 8.023911476135254
$ws.Range("E24").Value = 0.3867105368121521

$ws.Range("B25").Value = 2414.614165449143
$ws.Range("C25").Value = 50.81265128962131
$ws.Range("D25").Value = 8.284267520904541
$ws.Range("E25").Value = 0.6843681620468043

$ws.Range("B26").Value = 808.9159862995148
$ws.Range("C26").Value = 20.54685148138264
$ws.Range("D26").Value = 8.050402688980103
$ws.Range("E26").Value = 0.5464293443427336

$ws.Range("B27").Value = 1578.723355960846
$ws.Range("C27").Value = 33.44296924784374
$ws.Range("D27").Value = 7.555787897109985
$ws.Range("E27").Value = 0.9755899248595341

$ws.Range("B28").Value = 2109.100126647949
$ws.Range("C28").Value = 337.641018741166
$ws.Range("D28").Value = 2.320011043548584
$ws.Range("E28").Value = 1.514880365254806

